# Week 16 logging: a new player (N.Cottrell) recorded a stat for the
# first time this week, so a column for him needs to be inserted into
# both the "Rushing" and "Receiving" sheets, between T.Etienne and
# D.Chark (i.e. at column H), pushing the remaining player columns one
# place to the right. The new column gets the same placeholder data
# ("n") as every other player column in the "Yards list" row.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new blank column at H, shifting H:U (T.Etienne's
    # neighbours onward) one column to the right.
    $ws.Columns.Item(8).Insert()

    # Header for the newly inserted column.
    $ws.Cells.Item(1, 8).Value = "N.Cottrell"

    # Placeholder "Yards list" entry for the new player, matching the
    # rest of row 2.
    $ws.Cells.Item(2, 8).Value = "n"
}
